# Saldo_guide.xlsx update
# - Rename sheet (report re-generated one day later: 2024-10-21 -> 2024-10-22)
# - Bump every "Dt. Referencia" (column G) date by one day (45586 -> 45587)
# - For a subset of rows, the projected value (col D) was zeroed out and its
#   amount folded into the "Saldo Previsto" (col E), matching the already
#   unchanged "Vl. Total" (col H)
# - Row 108 had its Saldo Previsto / Vl. Total corrected from 95379.28 to 161.34
# - Selection moved to L15 on the single worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/table tab to the new export timestamp
$ws.Name = "IClientBalance-20241022-094147-"

# Every data row (2-274) has its reference date shifted from 45586 to 45587
$ws.Range("G2:G274").Value = 45587

# Rows where the projected balance (D) was moved entirely into the expected
# balance (E), so D becomes 0 and E becomes the (already-existing) total in H
$rowUpdates = @(
    @{ Row=5;   NewE=4354.3500000000004 },
    @{ Row=15;  NewE=12307.1 },
    @{ Row=17;  NewE=2994.85 },
    @{ Row=43;  NewE=3213.99 },
    @{ Row=49;  NewE=3243.6 },
    @{ Row=60;  NewE=5038.99 },
    @{ Row=104; NewE=27353.05 },
    @{ Row=132; NewE=2494.0700000000002 },
    @{ Row=143; NewE=17035.560000000001 },
    @{ Row=158; NewE=296.51 },
    @{ Row=173; NewE=7213.98 },
    @{ Row=235; NewE=2956.72 },
    @{ Row=264; NewE=8968.1 },
    @{ Row=265; NewE=5305.96 },
    @{ Row=270; NewE=3322.34 },
    @{ Row=271; NewE=4361.0200000000004 },
    @{ Row=273; NewE=3798.96 }
)

foreach ($item in $rowUpdates) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = $item.NewE
}

# Row 108 is a standalone correction: D108 was already 0, only E/H change
$ws.Cells.Item(108, 5).Value = 161.34
$ws.Cells.Item(108, 8).Value = 161.34

# Update the active cell selection to L15
$ws.Range("L15").Select()

Write-Host "edit applied"
